# daily auto push: 2026-01-06 13:47 UTC
# Insert two new rows of data right after the existing "2026/01/06" block
# (old row 567) and before the "2026/12/29" block (old row 568), shifting
# every subsequent row down by two. All other cell content stays the same,
# just renumbered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 568 (old row 568 and everything below
# shifts down by 2, ending up at 570.. / dimension grows to D611).
$ws.Rows.Item(568).Insert()
$ws.Rows.Item(568).Insert()

# The date strings in column A ("2026/01/06", "2026/12/29", ...) must stay
# plain text, exactly like every other row in the sheet -- otherwise Excel's
# COM layer auto-parses the "YYYY/MM/DD" pattern into a real date serial.
# Force text mode with the "@" number format while writing, then clear the
# formatting back off so the new cells end up styled exactly like their
# neighbours (no explicit style index), while keeping the literal string.
$ws.Range("A568:A569").NumberFormat = "@"

$ws.Range("A568").Value = "2026/01/06"
$ws.Range("B568").Value = "火"
$ws.Range("C568").Value = 18
$ws.Range("D568").Value = 17

$ws.Range("A569").Value = "2026/01/06"
$ws.Range("B569").Value = "火"
$ws.Range("C569").Value = 21
$ws.Range("D569").Value = 20

$ws.Range("A568:A569").ClearFormats()
